# Auto-generated edit script: restores market-data snapshot values (columns H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled-runner update.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 54
$ws.Cells.Item(5, 9).Value = 54
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 54
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = 61
$ws.Cells.Item(5, 14).Value = $null
$ws.Cells.Item(6, 8).Value = 1865.7222
$ws.Cells.Item(6, 9).Value = 398.85715
$ws.Cells.Item(6, 11).Value = 1196.57145
$ws.Cells.Item(6, 13).Value = -1084.57145
$ws.Cells.Item(9, 8).Value = 38.833332
$ws.Cells.Item(9, 9).Value = 48.25
$ws.Cells.Item(9, 11).Value = 48.25
$ws.Cells.Item(9, 13).Value = 120.75
$ws.Cells.Item(19, 8).Value = 39977.285
$ws.Cells.Item(19, 9).Value = 3349.6667
$ws.Cells.Item(19, 10).Value = 67448
$ws.Cells.Item(19, 11).Value = 3349.6667
$ws.Cells.Item(19, 12).Value = 67448
$ws.Cells.Item(19, 13).Value = -3174.6667
$ws.Cells.Item(19, 14).Value = -67798
$ws.Cells.Item(28, 8).Value = 49848.19
$ws.Cells.Item(28, 9).Value = 72087.07000000001
$ws.Cells.Item(28, 10).Value = 5370.4287
$ws.Cells.Item(28, 11).Value = 72087.07000000001
$ws.Cells.Item(28, 12).Value = 5370.4287
$ws.Cells.Item(28, 13).Value = -71602.07000000001
$ws.Cells.Item(28, 14).Value = -6340.4287
$ws.Cells.Item(63, 8).Value = 70499.5
$ws.Cells.Item(63, 9).Value = 65000
$ws.Cells.Item(63, 10).Value = 75999
$ws.Cells.Item(63, 11).Value = 65000
$ws.Cells.Item(63, 12).Value = 75999
$ws.Cells.Item(63, 13).Value = -64376
$ws.Cells.Item(63, 14).Value = -77247
$ws.Cells.Item(66, 8).Value = 70499.5
$ws.Cells.Item(66, 9).Value = 65000
$ws.Cells.Item(66, 10).Value = 75999
$ws.Cells.Item(66, 11).Value = 195000
$ws.Cells.Item(66, 12).Value = 227997
$ws.Cells.Item(66, 13).Value = -191880
$ws.Cells.Item(66, 14).Value = -234237
$ws.Cells.Item(80, 9).Value = 2284250.8
$ws.Cells.Item(80, 10).Value = 3038.2727
$ws.Cells.Item(80, 11).Value = 6852752.399999999
$ws.Cells.Item(80, 12).Value = 9114.8181
$ws.Cells.Item(80, 13).Value = -6851754.399999999
$ws.Cells.Item(80, 14).Value = -11110.8181
$ws.Cells.Item(83, 9).Value = 2284250.8
$ws.Cells.Item(83, 10).Value = 3038.2727
$ws.Cells.Item(83, 11).Value = 20558257.2
$ws.Cells.Item(83, 12).Value = 27344.4543
$ws.Cells.Item(83, 13).Value = -20553265.2
$ws.Cells.Item(83, 14).Value = -37328.4543
$ws.Cells.Item(86, 8).Value = 16753563
$ws.Cells.Item(86, 9).Value = 1988
$ws.Cells.Item(86, 10).Value = 25129350
$ws.Cells.Item(86, 11).Value = 1988
$ws.Cells.Item(86, 12).Value = 25129350
$ws.Cells.Item(86, 13).Value = -865
$ws.Cells.Item(86, 14).Value = -25131596
$ws.Cells.Item(88, 8).Value = 3186.6
$ws.Cells.Item(88, 10).Value = 3124.3333
$ws.Cells.Item(88, 12).Value = 3124.3333
$ws.Cells.Item(88, 14).Value = -3936.3333
$ws.Cells.Item(89, 8).Value = 16753563
$ws.Cells.Item(89, 9).Value = 1988
$ws.Cells.Item(89, 10).Value = 25129350
$ws.Cells.Item(89, 11).Value = 9940
$ws.Cells.Item(89, 12).Value = 125646750
$ws.Cells.Item(89, 13).Value = -4324
$ws.Cells.Item(89, 14).Value = -125657982
$ws.Cells.Item(91, 8).Value = 3186.6
$ws.Cells.Item(91, 10).Value = 3124.3333
$ws.Cells.Item(91, 12).Value = 3124.3333
$ws.Cells.Item(91, 14).Value = -5932.3333
$ws.Cells.Item(92, 8).Value = 91518.09
$ws.Cells.Item(92, 10).Value = 500750
$ws.Cells.Item(92, 12).Value = 500750
$ws.Cells.Item(92, 14).Value = -503246
$ws.Cells.Item(105, 8).Value = 29220.25
$ws.Cells.Item(105, 10).Value = 29220.25
$ws.Cells.Item(105, 12).Value = 29220.25
$ws.Cells.Item(105, 14).Value = -36208.25
$ws.Cells.Item(106, 8).Value = 127948.75
$ws.Cells.Item(106, 9).Value = 3480.6667
$ws.Cells.Item(106, 11).Value = 3480.6667
$ws.Cells.Item(106, 13).Value = -2849.6667
$ws.Cells.Item(107, 8).Value = 470.25
$ws.Cells.Item(107, 9).Value = 529.5
$ws.Cells.Item(107, 11).Value = 529.5
$ws.Cells.Item(107, 13).Value = 1390.5
$ws.Cells.Item(116, 8).Value = 41754690
$ws.Cells.Item(116, 9).Value = 22838618
$ws.Cells.Item(116, 10).Value = 111113620
$ws.Cells.Item(116, 11).Value = 22838618
$ws.Cells.Item(116, 12).Value = 111113620
$ws.Cells.Item(116, 13).Value = -22835176
$ws.Cells.Item(116, 14).Value = -111120504
$ws.Cells.Item(131, 8).Value = 10091.923
$ws.Cells.Item(131, 10).Value = 12010
$ws.Cells.Item(131, 12).Value = 36030
$ws.Cells.Item(131, 14).Value = -46110
$ws.Cells.Item(132, 8).Value = 3837.5
$ws.Cells.Item(132, 9).Value = 3524.9363
$ws.Cells.Item(132, 11).Value = 10574.8089
$ws.Cells.Item(132, 13).Value = -8044.8089
$ws.Cells.Item(135, 8).Value = 917.625
$ws.Cells.Item(135, 9).Value = 497.57144
$ws.Cells.Item(135, 10).Value = 1244.3334
$ws.Cells.Item(135, 11).Value = 4478.14296
$ws.Cells.Item(135, 12).Value = 11199.0006
$ws.Cells.Item(135, 13).Value = -1943.14296
$ws.Cells.Item(135, 14).Value = -16269.0006
$ws.Cells.Item(137, 8).Value = 1978.6522
$ws.Cells.Item(137, 9).Value = 1683.2354
$ws.Cells.Item(137, 11).Value = 5049.706200000001
$ws.Cells.Item(137, 13).Value = -2499.706200000001
$ws.Cells.Item(138, 8).Value = 2919.5
$ws.Cells.Item(138, 10).Value = 2831.3125
$ws.Cells.Item(138, 12).Value = 8493.9375
$ws.Cells.Item(138, 14).Value = -18773.9375
$ws.Cells.Item(141, 8).Value = 2200
$ws.Cells.Item(141, 9).Value = 1750
$ws.Cells.Item(141, 11).Value = 5250
$ws.Cells.Item(141, 13).Value = -70

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 7218.8887
$ws.Cells.Item(5, 9).Value = 1488.6666
$ws.Cells.Item(5, 10).Value = 10084
$ws.Cells.Item(5, 11).Value = 1488.6666
$ws.Cells.Item(5, 12).Value = 10084
$ws.Cells.Item(5, 13).Value = -1376.6666
$ws.Cells.Item(5, 14).Value = -10308
$ws.Cells.Item(32, 8).Value = 7468.8
$ws.Cells.Item(32, 9).Value = 2962.4656
$ws.Cells.Item(32, 10).Value = 19349.137
$ws.Cells.Item(32, 11).Value = 2962.4656
$ws.Cells.Item(32, 12).Value = 19349.137
$ws.Cells.Item(32, 13).Value = -2675.4656
$ws.Cells.Item(32, 14).Value = -19923.137
$ws.Cells.Item(34, 8).Value = 15000
$ws.Cells.Item(34, 9).Value = 15000
$ws.Cells.Item(34, 11).Value = 15000
$ws.Cells.Item(34, 13).Value = -14729
$ws.Cells.Item(45, 8).Value = 2068.1538
$ws.Cells.Item(45, 9).Value = 2219.111
$ws.Cells.Item(45, 10).Value = 1728.5
$ws.Cells.Item(45, 11).Value = 2219.111
$ws.Cells.Item(45, 12).Value = 1728.5
$ws.Cells.Item(45, 13).Value = -1842.111
$ws.Cells.Item(45, 14).Value = -2482.5
$ws.Cells.Item(50, 8).Value = 498.45456
$ws.Cells.Item(50, 9).Value = 120
$ws.Cells.Item(50, 11).Value = 120
$ws.Cells.Item(50, 13).Value = 594
$ws.Cells.Item(60, 8).Value = 35999
$ws.Cells.Item(60, 9).Value = 34999
$ws.Cells.Item(60, 11).Value = 34999
$ws.Cells.Item(60, 13).Value = -34266
$ws.Cells.Item(74, 8).Value = 16953596
$ws.Cells.Item(74, 9).Value = 20412494
$ws.Cells.Item(74, 10).Value = 5003.1
$ws.Cells.Item(74, 11).Value = 20412494
$ws.Cells.Item(74, 12).Value = 5003.1
$ws.Cells.Item(74, 13).Value = -20411620
$ws.Cells.Item(74, 14).Value = -6751.1
$ws.Cells.Item(77, 8).Value = 16953596
$ws.Cells.Item(77, 9).Value = 20412494
$ws.Cells.Item(77, 10).Value = 5003.1
$ws.Cells.Item(77, 11).Value = 102062470
$ws.Cells.Item(77, 12).Value = 25015.5
$ws.Cells.Item(77, 13).Value = -102058102
$ws.Cells.Item(77, 14).Value = -33751.5
$ws.Cells.Item(88, 8).Value = 8773146
$ws.Cells.Item(88, 9).Value = 18519280
$ws.Cells.Item(88, 10).Value = 1624.4
$ws.Cells.Item(88, 11).Value = 18519280
$ws.Cells.Item(88, 12).Value = 1624.4
$ws.Cells.Item(88, 13).Value = -18518874
$ws.Cells.Item(88, 14).Value = -2436.4
$ws.Cells.Item(91, 8).Value = 8773146
$ws.Cells.Item(91, 9).Value = 18519280
$ws.Cells.Item(91, 10).Value = 1624.4
$ws.Cells.Item(91, 11).Value = 18519280
$ws.Cells.Item(91, 12).Value = 1624.4
$ws.Cells.Item(91, 13).Value = -18517876
$ws.Cells.Item(91, 14).Value = -4432.4
$ws.Cells.Item(104, 8).Value = 37855.57
$ws.Cells.Item(104, 9).Value = 0
$ws.Cells.Item(104, 10).Value = 37855.57
$ws.Cells.Item(104, 11).Value = 0
$ws.Cells.Item(104, 12).Value = 37855.57
$ws.Cells.Item(104, 13).Value = $null
$ws.Cells.Item(104, 14).Value = -44843.57
$ws.Cells.Item(110, 8).Value = 32260638
$ws.Cells.Item(110, 9).Value = 43480920
$ws.Cells.Item(110, 11).Value = 43480920
$ws.Cells.Item(110, 13).Value = -43478875
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 13).Value = $null
$ws.Cells.Item(122, 8).Value = 4981.25
$ws.Cells.Item(122, 9).Value = 4962.8335
$ws.Cells.Item(122, 10).Value = 4999.6665
$ws.Cells.Item(122, 11).Value = 14888.5005
$ws.Cells.Item(122, 12).Value = 14998.9995
$ws.Cells.Item(122, 13).Value = -12438.5005
$ws.Cells.Item(122, 14).Value = -19898.9995
$ws.Cells.Item(132, 8).Value = 1311.2916
$ws.Cells.Item(132, 9).Value = 1064
$ws.Cells.Item(132, 10).Value = 6999
$ws.Cells.Item(132, 11).Value = 3192
$ws.Cells.Item(132, 12).Value = 20997
$ws.Cells.Item(132, 13).Value = -662
$ws.Cells.Item(132, 14).Value = -26057
$ws.Cells.Item(134, 8).Value = 61403.668
$ws.Cells.Item(134, 10).Value = 68500
$ws.Cells.Item(134, 12).Value = 68500
$ws.Cells.Item(134, 14).Value = -78640

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 7218.8887
$ws.Cells.Item(4, 9).Value = 1488.6666
$ws.Cells.Item(4, 10).Value = 10084
$ws.Cells.Item(4, 11).Value = 1488.6666
$ws.Cells.Item(4, 12).Value = 10084
$ws.Cells.Item(4, 13).Value = -1373.6666
$ws.Cells.Item(4, 14).Value = -10314
$ws.Cells.Item(20, 8).Value = 11318
$ws.Cells.Item(20, 9).Value = 12958.223
$ws.Cells.Item(20, 11).Value = 12958.223
$ws.Cells.Item(20, 13).Value = -12711.223
$ws.Cells.Item(24, 8).Value = 4950
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 13).Value = $null
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).Value = $null
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).Value = $null
$ws.Cells.Item(86, 8).Value = 1409.2
$ws.Cells.Item(86, 9).Value = 1143.85
$ws.Cells.Item(86, 10).Value = 2470.6
$ws.Cells.Item(86, 11).Value = 1143.85
$ws.Cells.Item(86, 12).Value = 2470.6
$ws.Cells.Item(86, 13).Value = -20.84999999999991
$ws.Cells.Item(86, 14).Value = -4716.6
$ws.Cells.Item(89, 8).Value = 1409.2
$ws.Cells.Item(89, 9).Value = 1143.85
$ws.Cells.Item(89, 10).Value = 2470.6
$ws.Cells.Item(89, 11).Value = 5719.25
$ws.Cells.Item(89, 12).Value = 12353
$ws.Cells.Item(89, 13).Value = -103.25
$ws.Cells.Item(89, 14).Value = -23585
$ws.Cells.Item(94, 8).Value = 16670320
$ws.Cells.Item(94, 9).Value = 27781060
$ws.Cells.Item(94, 10).Value = 4211.5
$ws.Cells.Item(94, 11).Value = 27781060
$ws.Cells.Item(94, 12).Value = 4211.5
$ws.Cells.Item(94, 13).Value = -27780609
$ws.Cells.Item(94, 14).Value = -5113.5
$ws.Cells.Item(105, 8).Value = 1649.303
$ws.Cells.Item(105, 9).Value = 1527.6818
$ws.Cells.Item(105, 10).Value = 1892.5454
$ws.Cells.Item(105, 11).Value = 1527.6818
$ws.Cells.Item(105, 12).Value = 1892.5454
$ws.Cells.Item(105, 13).Value = 219.3181999999999
$ws.Cells.Item(105, 14).Value = -5386.5454
$ws.Cells.Item(107, 8).Value = 62752170
$ws.Cells.Item(107, 9).Value = 503261.5
$ws.Cells.Item(107, 11).Value = 503261.5
$ws.Cells.Item(107, 13).Value = -501341.5
$ws.Cells.Item(134, 8).Value = 2717.3667
$ws.Cells.Item(134, 9).Value = 2552.4482
$ws.Cells.Item(134, 10).Value = 7500
$ws.Cells.Item(134, 11).Value = 7657.344599999999
$ws.Cells.Item(134, 12).Value = 22500
$ws.Cells.Item(134, 13).Value = -5122.344599999999
$ws.Cells.Item(134, 14).Value = -27570

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 27.142857
$ws.Cells.Item(7, 9).Value = 31.4
$ws.Cells.Item(7, 10).Value = 16.5
$ws.Cells.Item(7, 11).Value = 31.4
$ws.Cells.Item(7, 12).Value = 16.5
$ws.Cells.Item(7, 13).Value = 81.59999999999999
$ws.Cells.Item(7, 14).Value = -242.5
$ws.Cells.Item(16, 8).Value = 1063.9
$ws.Cells.Item(16, 10).Value = 1833
$ws.Cells.Item(16, 12).Value = 1833
$ws.Cells.Item(16, 14).Value = -2407
$ws.Cells.Item(18, 8).Value = 48699.668
$ws.Cells.Item(18, 10).Value = 48699.668
$ws.Cells.Item(18, 12).Value = 48699.668
$ws.Cells.Item(18, 14).Value = -49159.668
$ws.Cells.Item(31, 8).Value = 4971.9854
$ws.Cells.Item(31, 9).Value = 15540.889
$ws.Cells.Item(31, 10).Value = 3359.7795
$ws.Cells.Item(31, 11).Value = 15540.889
$ws.Cells.Item(31, 12).Value = 3359.7795
$ws.Cells.Item(31, 13).Value = -15245.889
$ws.Cells.Item(31, 14).Value = -3949.7795
$ws.Cells.Item(34, 8).Value = 4971.9854
$ws.Cells.Item(34, 9).Value = 15540.889
$ws.Cells.Item(34, 10).Value = 3359.7795
$ws.Cells.Item(34, 11).Value = 15540.889
$ws.Cells.Item(34, 12).Value = 3359.7795
$ws.Cells.Item(34, 13).Value = -15338.889
$ws.Cells.Item(34, 14).Value = -3763.7795
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).Value = $null
$ws.Cells.Item(88, 8).Value = 35879.8
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 35879.8
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 35879.8
$ws.Cells.Item(88, 13).Value = $null
$ws.Cells.Item(88, 14).Value = -36691.8
$ws.Cells.Item(91, 8).Value = 35879.8
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 35879.8
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 35879.8
$ws.Cells.Item(91, 13).Value = $null
$ws.Cells.Item(91, 14).Value = -38687.8
$ws.Cells.Item(99, 8).Value = 3083.1667
$ws.Cells.Item(99, 9).Value = 2875
$ws.Cells.Item(99, 10).Value = 3499.5
$ws.Cells.Item(99, 11).Value = 2875
$ws.Cells.Item(99, 12).Value = 3499.5
$ws.Cells.Item(99, 13).Value = -1377
$ws.Cells.Item(99, 14).Value = -6495.5
$ws.Cells.Item(113, 8).Value = 1063.9
$ws.Cells.Item(113, 10).Value = 1833
$ws.Cells.Item(113, 12).Value = 1833
$ws.Cells.Item(113, 14).Value = -6173
$ws.Cells.Item(122, 8).Value = 2627.4285
$ws.Cells.Item(122, 9).Value = 2918.8
$ws.Cells.Item(122, 10).Value = 1899
$ws.Cells.Item(122, 11).Value = 8756.400000000001
$ws.Cells.Item(122, 12).Value = 5697
$ws.Cells.Item(122, 13).Value = -6306.400000000001
$ws.Cells.Item(122, 14).Value = -10597
$ws.Cells.Item(126, 8).Value = 3083.1667
$ws.Cells.Item(126, 9).Value = 2875
$ws.Cells.Item(126, 10).Value = 3499.5
$ws.Cells.Item(126, 11).Value = 8625
$ws.Cells.Item(126, 12).Value = 10498.5
$ws.Cells.Item(126, 13).Value = -6155
$ws.Cells.Item(126, 14).Value = -15438.5
$ws.Cells.Item(132, 8).Value = 1252514.4
$ws.Cells.Item(132, 9).Value = 2186.1667
$ws.Cells.Item(132, 10).Value = 5003499
$ws.Cells.Item(132, 11).Value = 6558.500100000001
$ws.Cells.Item(132, 12).Value = 15010497
$ws.Cells.Item(132, 13).Value = -4028.500100000001
$ws.Cells.Item(132, 14).Value = -15015557
$ws.Cells.Item(134, 8).Value = 4347.727
$ws.Cells.Item(134, 10).Value = 4093
$ws.Cells.Item(134, 12).Value = 12279
$ws.Cells.Item(134, 14).Value = -17349
$ws.Cells.Item(141, 8).Value = 224998.9
$ws.Cells.Item(141, 10).Value = 224998.9
$ws.Cells.Item(141, 12).Value = 224998.9
$ws.Cells.Item(141, 14).Value = -235358.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 18828470
$ws.Cells.Item(4, 9).Value = 42358852
$ws.Cells.Item(4, 11).Value = 127076556
$ws.Cells.Item(4, 13).Value = -127076444
$ws.Cells.Item(70, 8).Value = 7500
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 7500
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 22500
$ws.Cells.Item(70, 13).Value = $null
$ws.Cells.Item(70, 14).Value = -23130
$ws.Cells.Item(73, 8).Value = 7500
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 7500
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 22500
$ws.Cells.Item(73, 13).Value = $null
$ws.Cells.Item(73, 14).Value = -24684
$ws.Cells.Item(75, 8).Value = 462.25
$ws.Cells.Item(75, 10).Value = 440.57144
$ws.Cells.Item(75, 12).Value = 1321.71432
$ws.Cells.Item(75, 14).Value = -3317.71432
$ws.Cells.Item(78, 8).Value = 462.25
$ws.Cells.Item(78, 10).Value = 440.57144
$ws.Cells.Item(78, 12).Value = 3965.14296
$ws.Cells.Item(78, 14).Value = -13949.14296
$ws.Cells.Item(103, 8).Value = 261.5
$ws.Cells.Item(103, 9).Value = 195.6
$ws.Cells.Item(103, 10).Value = 327.4
$ws.Cells.Item(103, 11).Value = 586.8
$ws.Cells.Item(103, 12).Value = 982.1999999999999
$ws.Cells.Item(103, 13).Value = 292.2
$ws.Cells.Item(103, 14).Value = -2740.2
$ws.Cells.Item(131, 8).Value = 2943.5588
$ws.Cells.Item(131, 9).Value = 1604.4286
$ws.Cells.Item(131, 10).Value = 3290.7407
$ws.Cells.Item(131, 11).Value = 4813.2858
$ws.Cells.Item(131, 12).Value = 9872.222099999999
$ws.Cells.Item(131, 13).Value = 226.7142000000003
$ws.Cells.Item(131, 14).Value = -19952.2221
$ws.Cells.Item(137, 8).Value = 1855.5
$ws.Cells.Item(137, 9).Value = 853.25
$ws.Cells.Item(137, 11).Value = 2559.75
$ws.Cells.Item(137, 13).Value = 2540.25
$ws.Cells.Item(139, 8).Value = 3716984.2
$ws.Cells.Item(139, 9).Value = 4777979.5
$ws.Cells.Item(139, 11).Value = 14333938.5
$ws.Cells.Item(139, 13).Value = -14328798.5
$ws.Cells.Item(141, 8).Value = 5377
$ws.Cells.Item(141, 9).Value = 5377
$ws.Cells.Item(141, 11).Value = 16131
$ws.Cells.Item(141, 13).Value = -10951

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 236.53847
$ws.Cells.Item(2, 9).Value = 42
$ws.Cells.Item(2, 11).Value = 42
$ws.Cells.Item(2, 13).Value = 71
$ws.Cells.Item(46, 8).Value = 9861
$ws.Cells.Item(46, 9).Value = 4149.8
$ws.Cells.Item(46, 10).Value = 17000
$ws.Cells.Item(46, 11).Value = 4149.8
$ws.Cells.Item(46, 12).Value = 17000
$ws.Cells.Item(46, 13).Value = -3993.8
$ws.Cells.Item(46, 14).Value = -17312
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).Value = $null
$ws.Cells.Item(80, 8).Value = 4355
$ws.Cells.Item(80, 9).Value = 4100
$ws.Cells.Item(80, 10).Value = 4567.5
$ws.Cells.Item(80, 11).Value = 4100
$ws.Cells.Item(80, 12).Value = 4567.5
$ws.Cells.Item(80, 13).Value = -3102
$ws.Cells.Item(80, 14).Value = -6563.5
$ws.Cells.Item(83, 8).Value = 4355
$ws.Cells.Item(83, 9).Value = 4100
$ws.Cells.Item(83, 10).Value = 4567.5
$ws.Cells.Item(83, 11).Value = 20500
$ws.Cells.Item(83, 12).Value = 22837.5
$ws.Cells.Item(83, 13).Value = -15508
$ws.Cells.Item(83, 14).Value = -32821.5
$ws.Cells.Item(102, 8).Value = 1886.381
$ws.Cells.Item(102, 9).Value = 1335.375
$ws.Cells.Item(102, 10).Value = 3649.6
$ws.Cells.Item(102, 11).Value = 1335.375
$ws.Cells.Item(102, 12).Value = 3649.6
$ws.Cells.Item(102, 13).Value = 286.625
$ws.Cells.Item(102, 14).Value = -6893.6
$ws.Cells.Item(113, 8).Value = 3474.5
$ws.Cells.Item(113, 9).Value = 949
$ws.Cells.Item(113, 11).Value = 949
$ws.Cells.Item(113, 13).Value = 1221
$ws.Cells.Item(122, 8).Value = 2470.9546
$ws.Cells.Item(122, 9).Value = 2288.9412
$ws.Cells.Item(122, 10).Value = 3089.8
$ws.Cells.Item(122, 11).Value = 6866.823600000001
$ws.Cells.Item(122, 12).Value = 9269.400000000001
$ws.Cells.Item(122, 13).Value = -4416.823600000001
$ws.Cells.Item(122, 14).Value = -14169.4
$ws.Cells.Item(126, 8).Value = 7221.75
$ws.Cells.Item(126, 9).Value = 8645.1875
$ws.Cells.Item(126, 11).Value = 25935.5625
$ws.Cells.Item(126, 13).Value = -23465.5625
$ws.Cells.Item(132, 8).Value = 5395.44
$ws.Cells.Item(132, 9).Value = 4624.3125
$ws.Cells.Item(132, 10).Value = 6766.3335
$ws.Cells.Item(132, 11).Value = 13872.9375
$ws.Cells.Item(132, 12).Value = 20299.0005
$ws.Cells.Item(132, 13).Value = -11342.9375
$ws.Cells.Item(132, 14).Value = -25359.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 31315488
$ws.Cells.Item(7, 9).Value = 50102380
$ws.Cells.Item(7, 10).Value = 4000
$ws.Cells.Item(7, 11).Value = 50102380
$ws.Cells.Item(7, 12).Value = 4000
$ws.Cells.Item(7, 13).Value = -50102268
$ws.Cells.Item(7, 14).Value = -4224
$ws.Cells.Item(16, 8).Value = 469
$ws.Cells.Item(16, 10).Value = 778.8
$ws.Cells.Item(16, 12).Value = 778.8
$ws.Cells.Item(16, 14).Value = -1118.8
$ws.Cells.Item(55, 8).Value = 850.5
$ws.Cells.Item(55, 9).Value = 1016.6667
$ws.Cells.Item(55, 10).Value = 352
$ws.Cells.Item(55, 11).Value = 1016.6667
$ws.Cells.Item(55, 12).Value = 352
$ws.Cells.Item(55, 13).Value = -843.6667
$ws.Cells.Item(55, 14).Value = -698
$ws.Cells.Item(61, 8).Value = 2490.125
$ws.Cells.Item(61, 9).Value = 2391.9644
$ws.Cells.Item(61, 10).Value = 3177.25
$ws.Cells.Item(61, 11).Value = 2391.9644
$ws.Cells.Item(61, 12).Value = 3177.25
$ws.Cells.Item(61, 13).Value = -2189.9644
$ws.Cells.Item(61, 14).Value = -3581.25
$ws.Cells.Item(82, 8).Value = 1171.8182
$ws.Cells.Item(82, 9).Value = 1156.6875
$ws.Cells.Item(82, 10).Value = 1212.1666
$ws.Cells.Item(82, 11).Value = 1156.6875
$ws.Cells.Item(82, 12).Value = 1212.1666
$ws.Cells.Item(82, 13).Value = -795.6875
$ws.Cells.Item(82, 14).Value = -1934.1666
$ws.Cells.Item(85, 8).Value = 1171.8182
$ws.Cells.Item(85, 9).Value = 1156.6875
$ws.Cells.Item(85, 10).Value = 1212.1666
$ws.Cells.Item(85, 11).Value = 1156.6875
$ws.Cells.Item(85, 12).Value = 1212.1666
$ws.Cells.Item(85, 13).Value = 91.3125
$ws.Cells.Item(85, 14).Value = -3708.1666
$ws.Cells.Item(92, 8).Value = 47176
$ws.Cells.Item(92, 9).Value = 25353
$ws.Cells.Item(92, 11).Value = 25353
$ws.Cells.Item(92, 13).Value = -22857
$ws.Cells.Item(113, 8).Value = 2490.125
$ws.Cells.Item(113, 9).Value = 2391.9644
$ws.Cells.Item(113, 10).Value = 3177.25
$ws.Cells.Item(113, 11).Value = 2391.9644
$ws.Cells.Item(113, 12).Value = 3177.25
$ws.Cells.Item(113, 13).Value = -221.9643999999998
$ws.Cells.Item(113, 14).Value = -7517.25
$ws.Cells.Item(126, 8).Value = 31315488
$ws.Cells.Item(126, 9).Value = 50102380
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 11).Value = 150307140
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = -150304670
$ws.Cells.Item(126, 14).Value = -16940
$ws.Cells.Item(136, 8).Value = 7249.8887
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 13).Value = -6450

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 694.0741
$ws.Cells.Item(107, 9).Value = 640.3182
$ws.Cells.Item(107, 10).Value = 930.6
$ws.Cells.Item(107, 11).Value = 1920.9546
$ws.Cells.Item(107, 12).Value = 2791.8
$ws.Cells.Item(107, 13).Value = -0.9546000000000276
$ws.Cells.Item(107, 14).Value = -6631.8
$ws.Cells.Item(113, 8).Value = 5690.613
$ws.Cells.Item(113, 10).Value = 3397.4736
$ws.Cells.Item(113, 12).Value = 10192.4208
$ws.Cells.Item(113, 14).Value = -14532.4208
$ws.Cells.Item(122, 8).Value = 2814.3
$ws.Cells.Item(122, 9).Value = 3207.818
$ws.Cells.Item(122, 11).Value = 9623.454000000002
$ws.Cells.Item(122, 13).Value = -7173.454000000002
$ws.Cells.Item(126, 8).Value = 1919.2
$ws.Cells.Item(126, 9).Value = 1749.125
$ws.Cells.Item(126, 10).Value = 2599.5
$ws.Cells.Item(126, 11).Value = 5247.375
$ws.Cells.Item(126, 12).Value = 7798.5
$ws.Cells.Item(126, 13).Value = -2777.375
$ws.Cells.Item(126, 14).Value = -12738.5
$ws.Cells.Item(132, 8).Value = 1006140.94
$ws.Cells.Item(132, 9).Value = 1181783.5
$ws.Cells.Item(132, 10).Value = 10833.333
$ws.Cells.Item(132, 11).Value = 3545350.5
$ws.Cells.Item(132, 12).Value = 32499.999
$ws.Cells.Item(132, 13).Value = -3542820.5
$ws.Cells.Item(132, 14).Value = -37559.999
$ws.Cells.Item(136, 8).Value = 4138.2
$ws.Cells.Item(136, 9).Value = 3987.6316
$ws.Cells.Item(136, 11).Value = 11962.8948
$ws.Cells.Item(136, 13).Value = -9412.8948

